$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook (after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Glucose and Kissat"

# Header row
$ws.Range("B2").Value = "Random Graphs"
$ws.Range("C2").Value = "Seed"
$ws.Range("D2").Value = "Glucose_4.1"
$ws.Range("E2").Value = "Kissat"

# Row 3
$ws.Range("B3").Value = "Edge = 130"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 72.062
$ws.Range("E3").Value = 9.51783

# Row 4
$ws.Range("B4").Value = "Density = 0.98"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 83.7909
$ws.Range("E4").Value = 13.6424

# Row 5
$ws.Range("B5").Value = "phase saving and variable elim disabled"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 107.567
$ws.Range("E5").Value = 11.6511

# Row 6
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 91.5744
$ws.Range("E6").Value = 18.5281

# Row 7
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 233.376
$ws.Range("E7").Value = 17.4116

# Row 8
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 80.2261
$ws.Range("E8").Value = 10.9032

# Row 9
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = 78.0786
$ws.Range("E9").Value = 8.78358

# Row 10
$ws.Range("C10").Value = 8
$ws.Range("D10").Value = 92.4235
$ws.Range("E10").Value = 10.2142

# Row 11
$ws.Range("C11").Value = 9
$ws.Range("D11").Value = 134.364
$ws.Range("E11").Value = 8.85007

# Row 12
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 110.367
$ws.Range("E12").Value = 13.7926

# Row 13 - Averages
$ws.Range("C13").Value = "Avgs"
$ws.Range("D13").Value = 108.38295
$ws.Range("E13").Value = 12.329468

# Column widths to fit the longer text labels (values chosen so the
# engine's pixel-rounded ColumnWidth lands as close as possible to the
# bestFit widths Excel computed: ~34.33 for column B, ~11.16 for column D)
$ws.Columns.Item(2).ColumnWidth = 33.42
$ws.Columns.Item(4).ColumnWidth = 10.3

# Selection / active cell on the new sheet
$ws.Range("E26").Select()
